try {
  "hello world" | Set-Content -Path "/tmp/work/psfile.txt"
  Write-Output "wrote"
} catch {
  Write-Output ("ERR1: " + $_.Exception.Message)
}
try {
  $x = Get-Content -Path "/tmp/work/psfile.txt"
  Write-Output ("read: " + $x)
} catch {
  Write-Output ("ERR2: " + $_.Exception.Message)
}
